{"js": "// Remove the empty paragraph, the \"Ver no Jupiter...\" paragraph, and the\n// \"\u00a9 2020 ...\" footer paragraph that used to follow the\n// \"LOT2059: Qu\u00edmica Org\u00e2nica Fundamental (Requisito fraco)\" requirement line.\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"text\");\nawait context.sync();\n\nconst targetTexts = [\n  \"Ver no Jupiter Salvar em pdf Salvar em docx\",\n  \"\u00a9 2020 . Contact: luizeleno@usp.br. Powered by Jekyll and Github pages. Original theme under Creative Commons Attribution\"\n];\n\n// Locate the \"Ver no Jupiter...\" paragraph; the empty paragraph right before\n// it (between it and the \"LOT2059...\" requirement paragraph) is removed too.\nlet verIndex = -1;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  if (paragraphs.items[i].text === targetTexts[0]) {\n    verIndex = i;\n    break;\n  }\n}\n\nif (verIndex === -1) {\n  throw new Error(\"Could not find the 'Ver no Jupiter...' paragraph\");\n}\n\nconst copyrightIndex = verIndex + 1;\nif (paragraphs.items[copyrightIndex].text !== targetTexts[1]) {\n  throw new Error(\"Unexpected paragraph after 'Ver no Jupiter...'\");\n}\n\nconst emptyIndex = verIndex - 1;\nif (paragraphs.items[emptyIndex].text !== \"\") {\n  throw new Error(\"Unexpected paragraph before 'Ver no Jupiter...'\");\n}\n\n// Delete from the end backwards so earlier indices stay valid.\nparagraphs.items[copyrightIndex].delete();\nparagraphs.items[verIndex].delete();\nparagraphs.items[emptyIndex].delete();\nawait context.sync();\n", "ps1": "# Remove the empty paragraph, the \"Ver no Jupiter...\" paragraph, and the\n# \"\u00a9 2020 ...\" footer paragraph that used to follow the\n# \"LOT2059: Qu\u00edmica Org\u00e2nica Fundamental (Requisito fraco)\" requirement line.\n$d = $word.ActiveDocument\n\n$verText = \"Ver no Jupiter Salvar em pdf Salvar em docx\"\n$copyrightText = \"\u00a9 2020 . Contact: luizeleno@usp.br. Powered by Jekyll and Github pages. Original theme under Creative Commons Attribution\"\n\n$verIndex = -1\nfor ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n    $txt = $d.Paragraphs.Item($i).Range.Text.TrimEnd([char]13, [char]7)\n    if ($txt -eq $verText) {\n        $verIndex = $i\n        break\n    }\n}\n\nif ($verIndex -eq -1) {\n    throw \"Could not find the 'Ver no Jupiter...' paragraph\"\n}\n\n$copyrightIndex = $verIndex + 1\n$copyrightActual = $d.Paragraphs.Item($copyrightIndex).Range.Text.TrimEnd([char]13, [char]7)\nif ($copyrightActual -ne $copyrightText) {\n    throw \"Unexpected paragraph after 'Ver no Jupiter...'\"\n}\n\n$emptyIndex = $verIndex - 1\n$emptyActual = $d.Paragraphs.Item($emptyIndex).Range.Text.TrimEnd([char]13, [char]7)\nif ($emptyActual -ne \"\") {\n    throw \"Unexpected paragraph before 'Ver no Jupiter...'\"\n}\n\n# Delete from the end backwards so earlier indices stay valid.\n$d.Paragraphs.Item($copyrightIndex).Range.Delete()\n$d.Paragraphs.Item($verIndex).Range.Delete()\n$d.Paragraphs.Item($emptyIndex).Range.Delete()\n"}
